# Scheduled-runner refresh of cached market price / leve-profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit* columns H:N) across the
# per-job Sheets. Values below are the new snapshot; a few rows also
# gain/lose a trailing LeveProfitHQ or LeveProfitNQ cell where the source
# feed newly has (or no longer has) an HQ/NQ price for that leve.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1739
$ws.Range("J2").Value = 2728.2856
$ws.Range("L2").Value = 2728.2856
$ws.Range("N2").Value = -2954.2856
$ws.Range("H8").Value = 33.857143
$ws.Range("I8").Value = 7.4
$ws.Range("K8").Value = 22.2
$ws.Range("M8").Value = 116.8
$ws.Range("H29").Value = 386.25
$ws.Range("I29").Value = 155.71428
$ws.Range("J29").Value = 2000
$ws.Range("K29").Value = 467.14284
$ws.Range("L29").Value = 6000
$ws.Range("M29").Value = -186.14284
$ws.Range("N29").Value = -6562
$ws.Range("H40").Value = 3037.9092
$ws.Range("J40").Value = 3181.8
$ws.Range("L40").Value = 3181.8
$ws.Range("N40").Value = -3531.8
$ws.Range("H70").Value = 752920.5
$ws.Range("I70").Value = 1559436.2
$ws.Range("J70").Value = 4013
$ws.Range("K70").Value = 4678308.6
$ws.Range("L70").Value = 12039
$ws.Range("M70").Value = -4678038.6
$ws.Range("N70").Value = -12579
$ws.Range("H73").Value = 752920.5
$ws.Range("I73").Value = 1559436.2
$ws.Range("J73").Value = 4013
$ws.Range("K73").Value = 4678308.6
$ws.Range("L73").Value = 12039
$ws.Range("M73").Value = -4677372.6
$ws.Range("N73").Value = -13911
$ws.Range("H92").Value = 1117
$ws.Range("I92").Value = 1161.6666
$ws.Range("J92").Value = 1050
$ws.Range("K92").Value = 1161.6666
$ws.Range("L92").Value = 1050
$ws.Range("M92").Value = 86.33339999999998
$ws.Range("N92").Value = -3546
$ws.Range("H98").Value = 4210
$ws.Range("J98").Value = 4079.1667
$ws.Range("L98").Value = 4079.1667
$ws.Range("N98").Value = -7075.1667
$ws.Range("H112").Value = 1872.9231
$ws.Range("I112").Value = 1219.6
$ws.Range("J112").Value = 2281.25
$ws.Range("K112").Value = 3658.8
$ws.Range("L112").Value = 6843.75
$ws.Range("M112").Value = -2550.8
$ws.Range("N112").Value = -9059.75
$ws.Range("H122").Value = 4210
$ws.Range("J122").Value = 4079.1667
$ws.Range("L122").Value = 12237.5001
$ws.Range("N122").Value = -17137.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7711.1055
$ws.Range("I32").Value = 7711.1055
$ws.Range("K32").Value = 7711.1055
$ws.Range("M32").Value = -7424.1055
$ws.Range("H45").Value = 2700.4
$ws.Range("J45").Value = 3497.5
$ws.Range("L45").Value = 3497.5
$ws.Range("N45").Value = -4251.5
$ws.Range("H60").Value = 2900
$ws.Range("I60").Value = 2900
$ws.Range("K60").Value = 2900
$ws.Range("M60").Value = -2167
$ws.Range("H61").Value = 2183.1667
$ws.Range("I61").Value = 1825
$ws.Range("J61").Value = 2899.5
$ws.Range("K61").Value = 1825
$ws.Range("L61").Value = 2899.5
$ws.Range("M61").Value = -1613
$ws.Range("N61").Value = -3323.5
$ws.Range("H97").Value = 1059.2667
$ws.Range("I97").Value = 840.9091
$ws.Range("J97").Value = 1659.75
$ws.Range("K97").Value = 840.9091
$ws.Range("L97").Value = 1659.75
$ws.Range("M97").Value = -344.9091
$ws.Range("N97").Value = -2651.75
$ws.Range("H122").Value = 3386.6875
$ws.Range("I122").Value = 3412.4666
$ws.Range("K122").Value = 10237.3998
$ws.Range("M122").Value = -7787.399800000001
$ws.Range("H125").Value = 33999
$ws.Range("J125").Value = 33999
$ws.Range("L125").Value = 33999
$ws.Range("N125").Value = -43839
$ws.Range("H136").Value = 2183.1667
$ws.Range("I136").Value = 1825
$ws.Range("J136").Value = 2899.5
$ws.Range("K136").Value = 5475
$ws.Range("L136").Value = 8698.5
$ws.Range("M136").Value = -2925
$ws.Range("N136").Value = -13798.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4330.0527
$ws.Range("I134").Value = 4215.0557
$ws.Range("J134").Value = 6400
$ws.Range("K134").Value = 12645.1671
$ws.Range("L134").Value = 19200
$ws.Range("M134").Value = -10110.1671
$ws.Range("N134").Value = -24270

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5887.875
$ws.Range("I62").Value = 5887.875
$ws.Range("K62").Value = 5887.875
$ws.Range("M62").Value = -5263.875
$ws.Range("H65").Value = 5887.875
$ws.Range("I65").Value = 5887.875
$ws.Range("K65").Value = 29439.375
$ws.Range("M65").Value = -26319.375
$ws.Range("H87").Value = 10000
$ws.Range("J87").Value = 10000
$ws.Range("L87").Value = 10000
$ws.Range("N87").Value = -12372
$ws.Range("H88").Value = 13098.7
$ws.Range("J88").Value = 13098.7
$ws.Range("L88").Value = 13098.7
$ws.Range("N88").Value = -13910.7
$ws.Range("H90").Value = 10000
$ws.Range("J90").Value = 10000
$ws.Range("L90").Value = 30000
$ws.Range("N90").Value = -41856
$ws.Range("H91").Value = 13098.7
$ws.Range("J91").Value = 13098.7
$ws.Range("L91").Value = 13098.7
$ws.Range("N91").Value = -15906.7
$ws.Range("H99").Value = 8000
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 8000
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 8000
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -10996
$ws.Range("H126").Value = 8000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 24000
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -28940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2222.25
$ws.Range("I70").Value = 2222.25
$ws.Range("K70").Value = 6666.75
$ws.Range("M70").Value = -6351.75
$ws.Range("H73").Value = 2222.25
$ws.Range("I73").Value = 2222.25
$ws.Range("K73").Value = 6666.75
$ws.Range("M73").Value = -5574.75
$ws.Range("H113").Value = 1522.1111
$ws.Range("I113").Value = 450
$ws.Range("J113").Value = 1828.4286
$ws.Range("K113").Value = 1350
$ws.Range("L113").Value = 5485.2858
$ws.Range("M113").Value = 820
$ws.Range("N113").Value = -9825.2858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3922.6667
$ws.Range("I102").Value = 3922.6667
$ws.Range("K102").Value = 3922.6667
$ws.Range("M102").Value = -2300.6667
$ws.Range("H113").Value = 1585.5
$ws.Range("I113").Value = 1585.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1585.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 584.5
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6326.909
$ws.Range("I40").Value = 4949.625
$ws.Range("J40").Value = 9999.666999999999
$ws.Range("K40").Value = 4949.625
$ws.Range("L40").Value = 9999.666999999999
$ws.Range("M40").Value = -4813.625
$ws.Range("N40").Value = -10271.667
$ws.Range("H64").Value = 29755
$ws.Range("J64").Value = 29755
$ws.Range("L64").Value = 29755
$ws.Range("N64").Value = -30205
$ws.Range("H67").Value = 29755
$ws.Range("J67").Value = 29755
$ws.Range("L67").Value = 29755
$ws.Range("N67").Value = -31315
$ws.Range("H82").Value = 1387.0714
$ws.Range("I82").Value = 1387.0714
$ws.Range("K82").Value = 1387.0714
$ws.Range("M82").Value = -1026.0714
$ws.Range("H85").Value = 1387.0714
$ws.Range("I85").Value = 1387.0714
$ws.Range("K85").Value = 1387.0714
$ws.Range("M85").Value = -139.0714
$ws.Range("H122").Value = 3200.7144
$ws.Range("I122").Value = 3200.7144
$ws.Range("K122").Value = 9602.143199999999
$ws.Range("M122").Value = -7152.143199999999
$ws.Range("H132").Value = 22846.55
$ws.Range("I132").Value = 22052
$ws.Range("K132").Value = 66156
$ws.Range("M132").Value = -63626

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 36994.5
$ws.Range("J31").Value = 36994.5
$ws.Range("L31").Value = 36994.5
$ws.Range("N31").Value = -37690.5
$ws.Range("H136").Value = 2500.25
$ws.Range("I136").Value = 2000.3334
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 6001.0002
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -3451.0002
$ws.Range("N136").Value = -17100
